$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# ---------------------------------------------------------------------
# Shape 2 ("TextBox 4") holds the main body copy on the title slide.
# ---------------------------------------------------------------------
$body = $s.Shapes.Item(2)
$tr = $body.TextFrame.TextRange

# --- Paragraph 1: re-type the trailing period so it lands in its own
#     run (purely cosmetic run split, no visible/formatting change). ---
$tailPeriod = $tr.Find("the functionality of the native level in the Web.")
$periodPos = $tailPeriod.Start + $tailPeriod.Length - 1
$period = $tr.Characters($periodPos, 1)
$period.Text = ""
$sentence = $tr.Find("the functionality of the native level in the Web")
$newPeriod = $sentence.InsertAfter(".")

# --- Paragraph 2: split the leading "A core feature of " into "A " +
#     "core feature of " (cosmetic split only). ---
$lead = $tr.Find("A core feature of ")
$lead.Text = ""
$afterLead = $tr.Find("the scheme is ")
$insCore = $afterLead.InsertBefore("core feature of ")
$afterLead2 = $tr.Find("core feature of ")
$insA = $afterLead2.InsertBefore("A ")

# --- Paragraph 2: re-type the trailing, italic period after
#     "innovation" so it becomes its own (still italic) run. ---
$tailPeriod2 = $tr.Find("innovation.")
$periodPos2 = $tailPeriod2.Start + $tailPeriod2.Length - 1
$period2 = $tr.Characters($periodPos2, 1)
$period2.Text = ""
$innovationWord = $tr.Find("innovation")
$newPeriod2 = $innovationWord.InsertAfter(".")
$newPeriod2.Font.Italic = $true

# ---------------------------------------------------------------------
# New paragraph inserted right after the "... innovation." paragraph:
# "In contrast to jumping between the Web and single-purpose "Apps",
#  a more seamless Web experience is facilitated."
# ("seamless Web " and "experience" are italic.)
# ---------------------------------------------------------------------
$para2 = $tr.Find("innovation.")
$quoteOpen = [char]0x201C
$quoteClose = [char]0x201D
$newParaText = "In contrast to jumping between the Web and single-purpose " + $quoteOpen + "Apps" + $quoteClose + ", a more seamless Web experience is facilitated."
$newPara = $para2.InsertAfter([char]13 + $newParaText)

$italicSpan1 = $tr.Find("seamless Web ")
$italicSpan1.Font.Italic = $false
$italicSpan1.Font.Italic = $true
$italicSpan2 = $tr.Find("experience is facilitated.")
$expWord = $tr.Characters($italicSpan2.Start, 10)
$expWord.Font.Italic = $true
$rest = $tr.Characters($italicSpan2.Start + 10, $italicSpan2.Length - 10)
$rest.Font.Italic = $false

$newParaLead = $tr.Find("In contrast to jumping between the Web and single-purpose ")
$newParaLead.Font.Italic = $false

# ---------------------------------------------------------------------
# Shape geometry: title/body text boxes shift upward and the body box
# grows taller to make room for the new paragraph.
# ---------------------------------------------------------------------
$title = $s.Shapes.Item(1)
$title.Top = 94.13489

$body.Top = 156.5041
$body.Height = 263.4515

$heading = $s.Shapes.Item(5)
$heading.Top = 48.87308
